$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 127
$ws.Range("H127").Value = 1595.75
$ws.Range("I127").Value = 601.0833
$ws.Range("J127").Value = 2590.4167
$ws.Range("K127").Value = 1803.2499
$ws.Range("L127").Value = 7771.250100000001
$ws.Range("M127").Value = 3156.7501
$ws.Range("N127").Value = -17691.2501
# Row 135
$ws.Range("H135").Value = 1315.7142
$ws.Range("I135").Value = 704.9167
$ws.Range("K135").Value = 6344.2503
$ws.Range("M135").Value = -3809.2503
# Row 137
$ws.Range("H137").Value = 1703387.4
$ws.Range("I137").Value = 3175837.8
$ws.Range("J137").Value = 4406.154
$ws.Range("K137").Value = 9527513.399999999
$ws.Range("L137").Value = 13218.462
$ws.Range("M137").Value = -9524963.399999999
$ws.Range("N137").Value = -18318.462
# Row 141
$ws.Range("H141").Value = 335299.66
$ws.Range("I141").Value = 501049.5
$ws.Range("J141").Value = 3800
$ws.Range("K141").Value = 1503148.5
$ws.Range("L141").Value = 11400
$ws.Range("M141").Value = -1497968.5
$ws.Range("N141").Value = -21760

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1778.8
$ws.Range("I61").Value = 1470
$ws.Range("J61").Value = 3014
$ws.Range("K61").Value = 1470
$ws.Range("L61").Value = 3014
$ws.Range("M61").Value = -1258
$ws.Range("N61").Value = -3438
# Row 74
$ws.Range("H74").Value = 7310.875
$ws.Range("I74").Value = 9039.4
$ws.Range("J74").Value = 4430
$ws.Range("K74").Value = 9039.4
$ws.Range("L74").Value = 4430
$ws.Range("M74").Value = -8165.4
$ws.Range("N74").Value = -6178
# Row 77
$ws.Range("H77").Value = 7310.875
$ws.Range("I77").Value = 9039.4
$ws.Range("J77").Value = 4430
$ws.Range("K77").Value = 45197
$ws.Range("L77").Value = 22150
$ws.Range("M77").Value = -40829
$ws.Range("N77").Value = -30886
# Row 122
$ws.Range("H122").Value = 3751.75
$ws.Range("I122").Value = 1700
$ws.Range("K122").Value = 5100
$ws.Range("M122").Value = -2650
# Row 132
$ws.Range("H132").Value = 3118.9333
$ws.Range("I132").Value = 1614
$ws.Range("J132").Value = 4838.857
$ws.Range("K132").Value = 4842
$ws.Range("L132").Value = 14516.571
$ws.Range("M132").Value = -2312
$ws.Range("N132").Value = -19576.571
# Row 136
$ws.Range("H136").Value = 1778.8
$ws.Range("I136").Value = 1470
$ws.Range("J136").Value = 3014
$ws.Range("K136").Value = 4410
$ws.Range("L136").Value = 9042
$ws.Range("M136").Value = -1860
$ws.Range("N136").Value = -14142

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3497.5789
$ws.Range("I134").Value = 2795
$ws.Range("J134").Value = 4702
$ws.Range("K134").Value = 8385
$ws.Range("L134").Value = 14106
$ws.Range("M134").Value = -5850
$ws.Range("N134").Value = -19176

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4194.7896
$ws.Range("J31").Value = 5538.4614
$ws.Range("L31").Value = 5538.4614
$ws.Range("N31").Value = -6128.4614
# Row 34
$ws.Range("H34").Value = 4194.7896
$ws.Range("J34").Value = 5538.4614
$ws.Range("L34").Value = 5538.4614
$ws.Range("N34").Value = -5942.4614
# Row 58
$ws.Range("H58").Value = 3068.8135
$ws.Range("I58").Value = 1833.8776
$ws.Range("J58").Value = 9120
$ws.Range("K58").Value = 1833.8776
$ws.Range("L58").Value = 9120
$ws.Range("M58").Value = -1630.8776
$ws.Range("N58").Value = -9526
# Row 107
$ws.Range("H107").Value = 660.7083
$ws.Range("I107").Value = 557.4
$ws.Range("K107").Value = 557.4
$ws.Range("M107").Value = 1362.6
# Row 132
$ws.Range("H132").Value = 4202.826
$ws.Range("J132").Value = 3994.9167
$ws.Range("L132").Value = 11984.7501
$ws.Range("N132").Value = -17044.7501
# Row 134
$ws.Range("H134").Value = 2558.4285
$ws.Range("I134").Value = 1281.8
$ws.Range("K134").Value = 3845.4
$ws.Range("M134").Value = -1310.4
# Row 136
$ws.Range("H136").Value = 3068.8135
$ws.Range("I136").Value = 1833.8776
$ws.Range("J136").Value = 9120
$ws.Range("K136").Value = 5501.6328
$ws.Range("L136").Value = 27360
$ws.Range("M136").Value = -2951.6328
$ws.Range("N136").Value = -32460

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 93.52941
$ws.Range("I12").Value = 12.5
$ws.Range("J12").Value = 118.46154
$ws.Range("K12").Value = 37.5
$ws.Range("L12").Value = 355.38462
$ws.Range("M12").Value = 135.5
$ws.Range("N12").Value = -701.38462
# Row 113
$ws.Range("H113").Value = 4465113.5
$ws.Range("I113").Value = 687.8570999999999
$ws.Range("J113").Value = 8929539
$ws.Range("K113").Value = 2063.5713
$ws.Range("L113").Value = 26788617
$ws.Range("M113").Value = 106.4287000000004
$ws.Range("N113").Value = -26792957

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 63
$ws.Range("H63").Value = 14900
$ws.Range("J63").Value = 14900
$ws.Range("L63").Value = 14900
$ws.Range("N63").Value = -16272
# Row 66
$ws.Range("H66").Value = 14900
$ws.Range("J66").Value = 14900
$ws.Range("L66").Value = 44700
$ws.Range("N66").Value = -51564
# Row 132
$ws.Range("H132").Value = 3692.2856
$ws.Range("I132").Value = 1521.5555
$ws.Range("J132").Value = 7599.6
$ws.Range("K132").Value = 4564.666499999999
$ws.Range("L132").Value = 22798.8
$ws.Range("M132").Value = -2034.666499999999
$ws.Range("N132").Value = -27858.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 15000000
$ws.Range("J2").Value = 15000000
$ws.Range("L2").Value = 15000000
$ws.Range("N2").Value = -15000224
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null
# Row 10
$ws.Range("H10").Value = 46499.5
$ws.Range("J10").Value = 46499.5
$ws.Range("L10").Value = 46499.5
$ws.Range("N10").Value = -46779.5
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null
# Row 20
$ws.Range("H20").Value = 10000
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10452
# Row 132
$ws.Range("H132").Value = 18031.572
$ws.Range("I132").Value = 26188
$ws.Range("J132").Value = 10616.637
$ws.Range("K132").Value = 78564
$ws.Range("L132").Value = 31849.911
$ws.Range("M132").Value = -76034
$ws.Range("N132").Value = -36909.911
# Row 136
$ws.Range("H136").Value = 5292.7
$ws.Range("I136").Value = 1728.2222
$ws.Range("J136").Value = 8209.091
$ws.Range("K136").Value = 5184.6666
$ws.Range("L136").Value = 24627.273
$ws.Range("M136").Value = -2634.6666
$ws.Range("N136").Value = -29727.273

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 13336488
$ws.Range("I132").Value = 2128.3125
$ws.Range("J132").Value = 37042016
$ws.Range("K132").Value = 6384.9375
$ws.Range("L132").Value = 111126048
$ws.Range("M132").Value = -3854.9375
$ws.Range("N132").Value = -111131108
# Row 136
$ws.Range("H136").Value = 8210.454
$ws.Range("J136").Value = 13500
$ws.Range("L136").Value = 40500
$ws.Range("N136").Value = -45600
